# Replace the "userdoc 'zone1'" field (fldChar/instrText run sequence) with
# plain-text runs spelling out the literal "{m:userdoc 'zone1'}" token, as
# produced by TokenIteratorFieldRewriterSplit.
$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph with the "userdoc 'zone1'" field -> literal "{m:userdoc 'zone1'}" ---
$p1 = $d.Paragraphs(2)
$body1 = '<w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ' + "'" + 'zone1' + "'" + '</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body>'
$p1.Range.InsertXML($xmlHeader + $body1 + $xmlFooter)

# --- Paragraph with the "enduserdoc" field -> literal "{m:enduserdoc}" (bookmark kept) ---
$p2 = $d.Paragraphs(4)
$body2 = '<w:body><w:p><w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">enduserdoc}</w:t></w:r></w:p></w:body>'
$p2.Range.InsertXML($xmlHeader + $body2 + $xmlFooter)

Write-Host "done"
